# Update "want to go" counts (column F) on several sheets, reflecting
# output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2718
$ws.Range("F4").Value = 1059
$ws.Range("F5").Value = 19465
$ws.Range("F7").Value = 2204
$ws.Range("F8").Value = 747
$ws.Range("F10").Value = 426
$ws.Range("F11").Value = 684
$ws.Range("F12").Value = 233
$ws.Range("F15").Value = 361
$ws.Range("F17").Value = 259
$ws.Range("F18").Value = 167
$ws.Range("F19").Value = 185
$ws.Range("F22").Value = 94

# --- Sheet: 演出 (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 190
$ws.Range("F5").Value = 14
$ws.Range("F7").Value = 283
$ws.Range("F8").Value = 126
$ws.Range("F21").Value = 33

# --- Sheet: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5993
$ws.Range("F3").Value = 637
$ws.Range("F4").Value = 584

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 5993
$ws.Range("F3").Value = 637
$ws.Range("F4").Value = 584
$ws.Range("F5").Value = 190
$ws.Range("F8").Value = 2718
$ws.Range("F9").Value = 1059
$ws.Range("F10").Value = 19466
$ws.Range("F12").Value = 14
$ws.Range("F15").Value = 283
$ws.Range("F16").Value = 2204
$ws.Range("F17").Value = 747
$ws.Range("F18").Value = 126
$ws.Range("F20").Value = 426
$ws.Range("F21").Value = 684
$ws.Range("F22").Value = 233
$ws.Range("F28").Value = 361
$ws.Range("F31").Value = 259
$ws.Range("F33").Value = 167
$ws.Range("F35").Value = 185
$ws.Range("F44").Value = 33
$ws.Range("F47").Value = 94
